$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the End Time (column D, row 3)
$ws.Range("D3").Value = "2:30pm"

# Update Hours Worked (column E, row 3) and its number format (0.0 -> 0.00)
$ws.Range("E2:E3").NumberFormat = "0.00"
$ws.Range("E3").Value = 4.75

# Update the Description (column F, row 3) with the revised/extended text
$ws.Range("F3").Value = "Started going through the codebase and adding comments to make the code more readable and understandable. Added comments to the follwing classes: EmailHandler, DataAcessor, AirQualitySensor. Restructured the folders in the visual studio solution. Was working in the wrong branch of git so I had to go back and re add them.  Created a timesheet to track hours worked. Still need to finish the AQS comments. Worked on the presentation by adding a theme, changing the fonts, and adding content to the road map"

# Increase row height to fit the longer description
$ws.Rows.Item(3).RowHeight = 75
